# Apply updated power-flow results for the 380 kV case (slack bus Vm set to 1.02 pu)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.026293730687815
$ws.Cells.Item(2, 4).Value = 1.03513027241954
$ws.Cells.Item(2, 5).Value = 1.029922842390774
$ws.Cells.Item(2, 6).Value = 1.041951930241646
$ws.Cells.Item(2, 9).Value = 1.031427613860127
$ws.Cells.Item(2, 10).Value = 1.031457977094533
$ws.Cells.Item(2, 11).Value = 1.037927759905831
$ws.Cells.Item(2, 12).Value = 1.03273534726391
$ws.Cells.Item(2, 13).Value = 1.04472999172796
$ws.Cells.Item(2, 14).Value = 1.014486559250125

# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.027254033722701
$ws.Cells.Item(3, 4).Value = 1.03587863893137
$ws.Cells.Item(3, 5).Value = 1.030830618003217
$ws.Cells.Item(3, 6).Value = 1.042890376771948
$ws.Cells.Item(3, 9).Value = 1.031571451937797
$ws.Cells.Item(3, 10).Value = 1.032057922865888
$ws.Cells.Item(3, 11).Value = 1.038485619794909
$ws.Cells.Item(3, 12).Value = 1.033451102157662
$ws.Cells.Item(3, 13).Value = 1.045478835792599
$ws.Cells.Item(3, 14).Value = 1.014688188467585

# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.027875906062445
$ws.Cells.Item(4, 4).Value = 1.036363090452392
$ws.Cells.Item(4, 5).Value = 1.03141885797848
$ws.Cells.Item(4, 6).Value = 1.043498324117879
$ws.Cells.Item(4, 9).Value = 1.03166321758378
$ws.Cells.Item(4, 10).Value = 1.032446028256833
$ws.Cells.Item(4, 11).Value = 1.038846137278419
$ws.Cells.Item(4, 12).Value = 1.033914463388416
$ws.Cells.Item(4, 13).Value = 1.045963473055198
$ws.Cells.Item(4, 14).Value = 1.014818542492957

# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.028137458166442
$ws.Cells.Item(5, 4).Value = 1.036566802404262
$ws.Cells.Item(5, 5).Value = 1.031666355974383
$ws.Cells.Item(5, 6).Value = 1.043754073766026
$ws.Cells.Item(5, 9).Value = 1.031701482390921
$ws.Cells.Item(5, 10).Value = 1.032609163022554
$ws.Cells.Item(5, 11).Value = 1.038997588842153
$ws.Cells.Item(5, 12).Value = 1.034109312092092
$ws.Cells.Item(5, 13).Value = 1.046167233771824
$ws.Cells.Item(5, 14).Value = 1.01487331577483

# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.028181380750933
$ws.Cells.Item(6, 4).Value = 1.036601009344182
$ws.Cells.Item(6, 5).Value = 1.03170792378189
$ws.Cells.Item(6, 6).Value = 1.043797025104897
$ws.Cells.Item(6, 9).Value = 1.031707888823198
$ws.Cells.Item(6, 10).Value = 1.032636552590999
$ws.Cells.Item(6, 11).Value = 1.039023011746881
$ws.Cells.Item(6, 12).Value = 1.034142031045668
$ws.Cells.Item(6, 13).Value = 1.046201447185706
$ws.Cells.Item(6, 14).Value = 1.014882510823729

# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.027879400475865
$ws.Cells.Item(7, 4).Value = 1.036365812272456
$ws.Cells.Item(7, 5).Value = 1.031422164266961
$ws.Cells.Item(7, 6).Value = 1.043501740796426
$ws.Cells.Item(7, 9).Value = 1.031663730112382
$ws.Cells.Item(7, 10).Value = 1.032448208169881
$ws.Cells.Item(7, 11).Value = 1.038848161415076
$ws.Cells.Item(7, 12).Value = 1.033917066764683
$ws.Cells.Item(7, 13).Value = 1.04596619564183
$ws.Cells.Item(7, 14).Value = 1.01481927448494

# Row 8
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.026618167210288
$ws.Cells.Item(8, 4).Value = 1.035383142437422
$ws.Cells.Item(8, 5).Value = 1.030229453002787
$ws.Cells.Item(8, 6).Value = 1.042268934902152
$ws.Cells.Item(8, 9).Value = 1.031476494963891
$ws.Cells.Item(8, 10).Value = 1.03166075167467
$ws.Cells.Item(8, 11).Value = 1.038116384900281
$ws.Cells.Item(8, 12).Value = 1.03297719382789
$ws.Cells.Item(8, 13).Value = 1.044983048773542
$ws.Cells.Item(8, 14).Value = 1.014554724083064

# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.024399518212295
$ws.Cells.Item(9, 4).Value = 1.033653207675718
$ws.Cells.Item(9, 5).Value = 1.028134294941488
$ws.Cells.Item(9, 6).Value = 1.040102066881865
$ws.Cells.Item(9, 9).Value = 1.031136574099597
$ws.Cells.Item(9, 10).Value = 1.030272428687651
$ws.Cells.Item(9, 11).Value = 1.036823460044744
$ws.Cells.Item(9, 12).Value = 1.031322746982049
$ws.Cells.Item(9, 13).Value = 1.043251321000644
$ws.Cells.Item(9, 14).Value = 1.01408769876192

# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.022923025052389
$ws.Cells.Item(10, 4).Value = 1.032501109885544
$ws.Cells.Item(10, 5).Value = 1.02674199889841
$ws.Cells.Item(10, 6).Value = 1.038661259053373
$ws.Cells.Item(10, 9).Value = 1.030903275170915
$ws.Cells.Item(10, 10).Value = 1.029346445141667
$ws.Cells.Item(10, 11).Value = 1.035959256338164
$ws.Cells.Item(10, 12).Value = 1.030221004637294
$ws.Cells.Item(10, 13).Value = 1.042097376131796
$ws.Cells.Item(10, 14).Value = 1.013775795909387

# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.022284314954369
$ws.Cells.Item(11, 4).Value = 1.032002538048196
$ws.Cells.Item(11, 5).Value = 1.026140195671061
$ws.Cells.Item(11, 6).Value = 1.03803828374254
$ws.Cells.Item(11, 9).Value = 1.030800674929123
$ws.Cells.Item(11, 10).Value = 1.028945391576342
$ws.Cells.Item(11, 11).Value = 1.035584524397617
$ws.Cells.Item(11, 12).Value = 1.02974424019707
$ws.Cells.Item(11, 13).Value = 1.041597847451733
$ws.Cells.Item(11, 14).Value = 1.013640611868497

# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.022047163622982
$ws.Cells.Item(12, 4).Value = 1.031817391755307
$ws.Cells.Item(12, 5).Value = 1.025916821085044
$ws.Cells.Item(12, 6).Value = 1.037807020040401
$ws.Cells.Item(12, 9).Value = 1.03076232769671
$ws.Cells.Item(12, 10).Value = 1.028796408739862
$ws.Cells.Item(12, 11).Value = 1.035445254200457
$ws.Cells.Item(12, 12).Value = 1.029567194269394
$ws.Cells.Item(12, 13).Value = 1.041412321943652
$ws.Cells.Item(12, 14).Value = 1.013590379602501

# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.022098029121711
$ws.Cells.Item(13, 4).Value = 1.031857104184404
$ws.Cells.Item(13, 5).Value = 1.025964728345941
$ws.Cells.Item(13, 6).Value = 1.037856620666889
$ws.Cells.Item(13, 9).Value = 1.03077056402225
$ws.Cells.Item(13, 10).Value = 1.028828366670061
$ws.Cells.Item(13, 11).Value = 1.035475131647885
$ws.Cells.Item(13, 12).Value = 1.029605169142466
$ws.Cells.Item(13, 13).Value = 1.041452116796052
$ws.Cells.Item(13, 14).Value = 1.013601155443559

# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.022264710028861
$ws.Cells.Item(14, 4).Value = 1.031987232862444
$ws.Cells.Item(14, 5).Value = 1.026121728148159
$ws.Cells.Item(14, 6).Value = 1.038019164609418
$ws.Cells.Item(14, 9).Value = 1.030797509965521
$ws.Cells.Item(14, 10).Value = 1.028933076883006
$ws.Cells.Item(14, 11).Value = 1.035573013866708
$ws.Cells.Item(14, 12).Value = 1.029729604589444
$ws.Cells.Item(14, 13).Value = 1.041582511406661
$ws.Cells.Item(14, 14).Value = 1.013636460037167

# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.022367420119201
$ws.Cells.Item(15, 4).Value = 1.032067415502036
$ws.Cells.Item(15, 5).Value = 1.026218482399413
$ws.Cells.Item(15, 6).Value = 1.038119331497142
$ws.Cells.Item(15, 9).Value = 1.030814080868861
$ws.Cells.Item(15, 10).Value = 1.028997590512315
$ws.Cells.Item(15, 11).Value = 1.035633312016035
$ws.Cells.Item(15, 12).Value = 1.02980627944795
$ws.Cells.Item(15, 13).Value = 1.041662854736041
$ws.Cells.Item(15, 14).Value = 1.013658209869107

# Row 16
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.022965427286312
$ws.Cells.Item(16, 4).Value = 1.032534204763887
$ws.Cells.Item(16, 5).Value = 1.026781961241548
$ws.Cells.Item(16, 6).Value = 1.038702623025873
$ws.Cells.Item(16, 9).Value = 1.03091005115785
$ws.Cells.Item(16, 10).Value = 1.029373059781496
$ws.Cells.Item(16, 11).Value = 1.035984115062017
$ws.Cells.Item(16, 12).Value = 1.030252652267325
$ws.Cells.Item(16, 13).Value = 1.042130531186167
$ws.Cells.Item(16, 14).Value = 1.013784764960476

# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.023340708192487
$ws.Cells.Item(17, 4).Value = 1.032827088994575
$ws.Cells.Item(17, 5).Value = 1.02713570395274
$ws.Cells.Item(17, 6).Value = 1.039068749298188
$ws.Cells.Item(17, 9).Value = 1.03096982804301
$ws.Cells.Item(17, 10).Value = 1.029608556458962
$ws.Cells.Item(17, 11).Value = 1.036204024389023
$ws.Cells.Item(17, 12).Value = 1.030532730402369
$ws.Cells.Item(17, 13).Value = 1.042423929734443
$ws.Cells.Item(17, 14).Value = 1.013864115586359

# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.023559663083009
$ws.Cells.Item(18, 4).Value = 1.032997951701193
$ws.Cells.Item(18, 5).Value = 1.027342139267772
$ws.Cells.Item(18, 6).Value = 1.039282391696009
$ws.Cells.Item(18, 9).Value = 1.031004542383575
$ws.Cells.Item(18, 10).Value = 1.029745908312766
$ws.Cells.Item(18, 11).Value = 1.036332242959885
$ws.Cells.Item(18, 12).Value = 1.030696123857839
$ws.Cells.Item(18, 13).Value = 1.042595077294969
$ws.Cells.Item(18, 14).Value = 1.013910387086171

# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.023634331178705
$ws.Cells.Item(19, 4).Value = 1.033056216227589
$ws.Cells.Item(19, 5).Value = 1.027412545840969
$ws.Cells.Item(19, 6).Value = 1.039355252974608
$ws.Cells.Item(19, 9).Value = 1.031016353194409
$ws.Cells.Item(19, 10).Value = 1.029792740130386
$ws.Cells.Item(19, 11).Value = 1.036375953530073
$ws.Cells.Item(19, 12).Value = 1.030751841628813
$ws.Cells.Item(19, 13).Value = 1.0426534363877
$ws.Cells.Item(19, 14).Value = 1.013926162361908

# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.023300437903275
$ws.Cells.Item(20, 4).Value = 1.032795662352862
$ws.Cells.Item(20, 5).Value = 1.027097740027254
$ws.Cells.Item(20, 6).Value = 1.03902945839203
$ws.Cells.Item(20, 9).Value = 1.030963430319788
$ws.Cells.Item(20, 10).Value = 1.02958329085779
$ws.Cells.Item(20, 11).Value = 1.036180435448919
$ws.Cells.Item(20, 12).Value = 1.030502677702125
$ws.Cells.Item(20, 13).Value = 1.042392449478692
$ws.Cells.Item(20, 14).Value = 1.013855603292969

# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.022215624054478
$ws.Cells.Item(21, 4).Value = 1.031948911944049
$ws.Cells.Item(21, 5).Value = 1.026075491140273
$ws.Cells.Item(21, 6).Value = 1.037971295675878
$ws.Cells.Item(21, 9).Value = 1.030789581593502
$ws.Cells.Item(21, 10).Value = 1.028902242706283
$ws.Cells.Item(21, 11).Value = 1.035544192140295
$ws.Cells.Item(21, 12).Value = 1.029692960184972
$ws.Cells.Item(21, 13).Value = 1.041544112836435
$ws.Cells.Item(21, 14).Value = 1.013626064233106

# Row 22
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.021534102960369
$ws.Cells.Item(22, 4).Value = 1.031416790167991
$ws.Cells.Item(22, 5).Value = 1.025433700127728
$ws.Cells.Item(22, 6).Value = 1.037306780422268
$ws.Cells.Item(22, 9).Value = 1.030678905410057
$ws.Cells.Item(22, 10).Value = 1.028473961932363
$ws.Cells.Item(22, 11).Value = 1.035143709509848
$ws.Cells.Item(22, 12).Value = 1.029184123021865
$ws.Cells.Item(22, 13).Value = 1.04101085592502
$ws.Cells.Item(22, 14).Value = 1.01348163477016

# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.021895338342541
$ws.Cells.Item(23, 4).Value = 1.031698852549215
$ws.Cells.Item(23, 5).Value = 1.025773836443942
$ws.Cells.Item(23, 6).Value = 1.037658976884317
$ws.Cells.Item(23, 9).Value = 1.03073770670349
$ws.Cells.Item(23, 10).Value = 1.028701008864797
$ws.Cells.Item(23, 11).Value = 1.035356055378086
$ws.Cells.Item(23, 12).Value = 1.029453841851142
$ws.Cells.Item(23, 13).Value = 1.04129353319411
$ws.Cells.Item(23, 14).Value = 1.013558209792208

# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.023318634128527
$ws.Cells.Item(24, 4).Value = 1.03280986261222
$ws.Cells.Item(24, 5).Value = 1.027114893972342
$ws.Cells.Item(24, 6).Value = 1.039047211992129
$ws.Cells.Item(24, 9).Value = 1.03096632164678
$ws.Cells.Item(24, 10).Value = 1.029594707324437
$ws.Cells.Item(24, 11).Value = 1.036191094432678
$ws.Cells.Item(24, 12).Value = 1.030516257135411
$ws.Cells.Item(24, 13).Value = 1.042406674010666
$ws.Cells.Item(24, 14).Value = 1.01385944967022

# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.02497263632027
$ws.Cells.Item(25, 4).Value = 1.034100232973267
$ws.Cells.Item(25, 5).Value = 1.028675159805579
$ws.Cells.Item(25, 6).Value = 1.0406615952339
$ws.Cells.Item(25, 9).Value = 1.03122563232665
$ws.Cells.Item(25, 10).Value = 1.030631423914184
$ws.Cells.Item(25, 11).Value = 1.037158113740643
$ws.Cells.Item(25, 12).Value = 1.031750250084781
$ws.Cells.Item(25, 13).Value = 1.043698923880785
$ws.Cells.Item(25, 14).Value = 1.01420853497087
